# Scope of the Project.docx — add an "Introduction" heading.
#
# The document currently has an empty, centered, 20pt bold paragraph
# directly under the title ("Scope Management Plan"). We turn that
# paragraph into a left-aligned, 16pt bold-underlined "Introduction "
# heading, and insert a new empty paragraph (with the same 16pt
# bold-underlined formatting, but no centering) right after it — taking
# over the spot the original blank paragraph used to occupy.

$d = $word.ActiveDocument
$target = $d.Paragraphs.Item(2)

# 1) Drop the centered alignment (paragraph becomes left/default aligned).
$target.Alignment = 0

# 2) Re-point the run-level formatting (size 16pt, single underline) that
#    lives on the paragraph mark. The COM shim only propagates Font
#    writes onto a Range that currently owns real text, so stage a
#    throw-away character, restyle it, then delete just that character
#    (not the whole paragraph) so the paragraph goes back to being truly
#    empty while keeping the new formatting on its mark.
$target.Range.Text = "X"
$target.Range.Font.Size = 16
$target.Range.Font.Underline = 1
$d.Range($target.Range.Start, $target.Range.Start + 1).Delete()

# 3) Insert a fresh paragraph right before this (now reformatted, still
#    empty) paragraph; it inherits the same pPr/rPr, and becomes the new
#    "Introduction" heading, while $target itself slides down to become
#    the trailing blank paragraph called for by the diff.
$target.Range.InsertParagraphBefore()

$intro = $d.Paragraphs.Item(2)
$intro.Range.Text = "Introduction "
